# Buff_buff表.xlsx — "feat: flying first buff"
#
# Adds a new buff-table column `flyFirst` (AI) right after the existing
# `multiHit` column (AH), and adds a new buff row (id 4011, "优先对空")
# right after buff 4010, before the 5001/5002/5003 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert the new data row (old row 35 "5001" .. row 37 "5003" shift
#    down to 36..38) and populate it as buff id 4011 "优先对空".
# ---------------------------------------------------------------------
$ws.Rows.Item(35).Insert()

$ws.Range("A35").Value = 4011
$ws.Range("B35").Value = 4011
$ws.Range("D35").Value = "优先对空"
$ws.Range("E35").Value = -1
$ws.Range("F35").Value = 1
$ws.Range("H35").Value = 1

$zeroCols35 = @("I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH")
foreach ($col in $zeroCols35) {
    $ws.Range($col + "35").Value = 0
}

# ---------------------------------------------------------------------
# 2) Insert the new `flyFirst` column AI (right after `multiHit` / AH)
#    and fill in its header rows.
# ---------------------------------------------------------------------
$ws.Columns.Item(35).Insert()

$ws.Range("AI1").Value = "int"
$ws.Range("AI2").Value = "flyFirst"
$ws.Range("AI3").Value = "优先对空`n0 没有`n1 有"

# ---------------------------------------------------------------------
# 3) Fill in the flyFirst (AI) value for every data row. Default is 0;
#    the new 4011 "优先对空" row (35) is the only one turned on (1).
# ---------------------------------------------------------------------
for ($r = 5; $r -le 38; $r++) {
    if ($r -eq 35) {
        $ws.Range("AI" + $r).Value = 1
    } else {
        $ws.Range("AI" + $r).Value = 0
    }
}

# ---------------------------------------------------------------------
# 4) View-state touch-ups captured by the diff: zoom 80% -> 60%, and the
#    remembered selection moves from AH2 to AG49.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 60
$ws.Range("AG49").Select()
